$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.483.30"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "2.039.13"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'230.51"
$ws.Range("E5").Value = "  +2.16%  "
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'56.35"
$ws.Range("E8").Value = "  +3.02%  "
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").Value = "'0.0801"
$ws.Range("E10").Value = "  +1.76%  "
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("D12").Value = "2.336.15"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "'14.48"
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("D14").Value = "'20.40"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("D15").Value = "'0.746"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("E16").Value = "  +1.89%  "
$ws.Range("D17").Value = "2.037.17"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").Value = "37.405.87"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("D19").Value = "'6.24"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").Value = "'69.12"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("D21").Value = "0.0₃0826"
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("D22").Value = "'223.75"
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  +2.01%  "
$ws.Range("E25").Value = "  +2.40%  "
$ws.Range("D26").Value = "'164.76"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("E28").Value = "  +6.74%  "
$ws.Range("D29").Value = "'18.81"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("E32").Value = "  +0.67%  "
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("E34").Value = "  +10.48%  "
$ws.Range("D35").Value = "'4.50"
$ws.Range("E35").Value = "  +1.55%  "
$ws.Range("D36").Value = "'2.35"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D37").Value = "'5.75"
$ws.Range("E37").Value = "  +8.90%  "
$ws.Range("D38").Value = "'3.23"
$ws.Range("E38").Value = "  +2.47%  "
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").Value = "1.473.71"
$ws.Range("E40").Value = "  -1.09%  "
$ws.Range("D41").Value = "'0.0214"
$ws.Range("E41").Value = "  -1.80%  "
$ws.Range("D42").Value = "'0.0933"
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("D43").Value = "'94.91"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  +2.43%  "
$ws.Range("D45").Value = "'4.24"
$ws.Range("E45").Value = "  +17.14%  "
$ws.Range("E46").Value = "  -4.46%  "
$ws.Range("E47").Value = "  -1.56%  "
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("E49").Value = "  -2.69%  "
$ws.Range("D50").Value = "'2.95"
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("D51").Value = "2.223.21"
$ws.Range("E51").Value = "  +0.74%  "

# Reset style to Normal for cells forced to text via apostrophe prefix,
# to avoid leaving an explicit text-format style applied.
foreach ($addr in @("D5","D8","D10","D13","D14","D15","D19","D20","D22","D23","D26","D29","D35","D36","D37","D38","D41","D42","D43","D45","D50")) {
    $ws.Range($addr).Style = "Normal"
}
